$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33500
$ws.Range("J3").Value = 33500
$ws.Range("L3").Value = 33500
$ws.Range("N3").Value = -33728
$ws.Range("H17").Value = 3340583.5
$ws.Range("J17").Value = 3340583.5
$ws.Range("L17").Value = 10021750.5
$ws.Range("N17").Value = -10022086.5
$ws.Range("H19").Value = 845.1667
$ws.Range("I19").Value = 950
$ws.Range("J19").Value = 824.2
$ws.Range("K19").Value = 950
$ws.Range("L19").Value = 824.2
$ws.Range("M19").Value = -775
$ws.Range("N19").Value = -1174.2
$ws.Range("H51").Value = 6590.5
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968
$ws.Range("H70").Value = 3713.9
$ws.Range("I70").Value = 4096.25
$ws.Range("K70").Value = 12288.75
$ws.Range("M70").Value = -12018.75
$ws.Range("H73").Value = 3713.9
$ws.Range("I73").Value = 4096.25
$ws.Range("K73").Value = 12288.75
$ws.Range("M73").Value = -11352.75
$ws.Range("H87").Value = 354354
$ws.Range("J87").Value = 354354
$ws.Range("L87").Value = 354354
$ws.Range("N87").Value = -356850
$ws.Range("H90").Value = 354354
$ws.Range("J90").Value = 354354
$ws.Range("L90").Value = 1063062
$ws.Range("N90").Value = -1075542
$ws.Range("H102").Value = 33500
$ws.Range("J102").Value = 33500
$ws.Range("L102").Value = 33500
$ws.Range("N102").Value = -39990
$ws.Range("H106").Value = 7238.4
$ws.Range("I106").Value = 8224.5
$ws.Range("J106").Value = 3294
$ws.Range("K106").Value = 8224.5
$ws.Range("L106").Value = 3294
$ws.Range("M106").Value = -7593.5
$ws.Range("N106").Value = -4556
$ws.Range("H137").Value = 9325.925999999999
$ws.Range("I137").Value = 1274.7
$ws.Range("K137").Value = 3824.1
$ws.Range("M137").Value = -1274.1
$ws.Range("H138").Value = 2268.8682
$ws.Range("I138").Value = 1562.5555
$ws.Range("J138").Value = 2566.8438
$ws.Range("K138").Value = 4687.666499999999
$ws.Range("L138").Value = 7700.5314
$ws.Range("M138").Value = 452.3335000000006
$ws.Range("N138").Value = -17980.5314
$ws.Range("H141").Value = 3035.64
$ws.Range("I141").Value = 3295.1667
$ws.Range("K141").Value = 9885.500100000001
$ws.Range("M141").Value = -4705.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2165.8
$ws.Range("I2").Value = 1172.3636
$ws.Range("J2").Value = 2946.3572
$ws.Range("K2").Value = 1172.3636
$ws.Range("L2").Value = 2946.3572
$ws.Range("M2").Value = -1059.3636
$ws.Range("N2").Value = -3172.3572
$ws.Range("H116").Value = 2165.8
$ws.Range("I116").Value = 1172.3636
$ws.Range("J116").Value = 2946.3572
$ws.Range("K116").Value = 1172.3636
$ws.Range("L116").Value = 2946.3572
$ws.Range("M116").Value = 1121.6364
$ws.Range("N116").Value = -7534.3572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2165.8
$ws.Range("I3").Value = 1172.3636
$ws.Range("J3").Value = 2946.3572
$ws.Range("K3").Value = 1172.3636
$ws.Range("L3").Value = 2946.3572
$ws.Range("M3").Value = -1058.3636
$ws.Range("N3").Value = -3174.3572
$ws.Range("H94").Value = 2089.4524
$ws.Range("I94").Value = 1086.7059
$ws.Range("J94").Value = 6351.125
$ws.Range("K94").Value = 1086.7059
$ws.Range("L94").Value = 6351.125
$ws.Range("M94").Value = -635.7058999999999
$ws.Range("N94").Value = -7253.125
$ws.Range("H105").Value = 58824920
$ws.Range("I105").Value = 66668070
$ws.Range("K105").Value = 66668070
$ws.Range("M105").Value = -66666323
$ws.Range("H134").Value = 36084.8
$ws.Range("I134").Value = 45763.566
$ws.Range("K134").Value = 137290.698
$ws.Range("M134").Value = -134755.698

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10704.442
$ws.Range("I58").Value = 3872.9167
$ws.Range("K58").Value = 3872.9167
$ws.Range("M58").Value = -3669.9167
$ws.Range("H122").Value = 2106.2856
$ws.Range("I122").Value = 1776.8235
$ws.Range("J122").Value = 3506.5
$ws.Range("K122").Value = 5330.470499999999
$ws.Range("L122").Value = 10519.5
$ws.Range("M122").Value = -2880.470499999999
$ws.Range("N122").Value = -15419.5
$ws.Range("H136").Value = 10704.442
$ws.Range("I136").Value = 3872.9167
$ws.Range("K136").Value = 11618.7501
$ws.Range("M136").Value = -9068.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9292812
$ws.Range("I5").Value = 1253.8182
$ws.Range("J5").Value = 23893832
$ws.Range("K5").Value = 3761.4546
$ws.Range("L5").Value = 71681496
$ws.Range("M5").Value = -3649.4546
$ws.Range("N5").Value = -71681720
$ws.Range("H16").Value = 9.5
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H119").Value = 7404
$ws.Range("H131").Value = 1408.34
$ws.Range("I131").Value = 341.66666
$ws.Range("J131").Value = 1441.33
$ws.Range("K131").Value = 1024.99998
$ws.Range("L131").Value = 4323.99
$ws.Range("M131").Value = 4015.00002
$ws.Range("N131").Value = -14403.99
$ws.Range("H135").Value = 9292812
$ws.Range("I135").Value = 1253.8182
$ws.Range("J135").Value = 23893832
$ws.Range("K135").Value = 11284.3638
$ws.Range("L135").Value = 215044488
$ws.Range("M135").Value = -8749.363799999999
$ws.Range("N135").Value = -215049558
$ws.Range("H139").Value = 15005.5
$ws.Range("I139").Value = 19883.375
$ws.Range("J139").Value = 5249.75
$ws.Range("K139").Value = 59650.125
$ws.Range("L139").Value = 15749.25
$ws.Range("M139").Value = -54510.125
$ws.Range("N139").Value = -26029.25
$ws.Range("N16").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 7949
$ws.Range("I20").Value = 2156.5
$ws.Range("K20").Value = 2156.5
$ws.Range("M20").Value = -1911.5
$ws.Range("H24").Value = 12060.625
$ws.Range("I24").Value = 7161.6665
$ws.Range("J24").Value = 15000
$ws.Range("K24").Value = 7161.6665
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = -6988.6665
$ws.Range("N24").Value = -15346
$ws.Range("H102").Value = 3976522.8
$ws.Range("I102").Value = 6436994
$ws.Range("K102").Value = 6436994
$ws.Range("M102").Value = -6435372
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080
$ws.Range("H126").Value = 5657487.5
$ws.Range("I126").Value = 4500747
$ws.Range("J126").Value = 6498753.5
$ws.Range("K126").Value = 13502241
$ws.Range("L126").Value = 19496260.5
$ws.Range("M126").Value = -13499771
$ws.Range("N126").Value = -19501200.5
$ws.Range("H132").Value = 4710
$ws.Range("I132").Value = 4710
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14130
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -11600
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 39999.332
$ws.Range("I62").Value = 40000
$ws.Range("J62").Value = 39999
$ws.Range("K62").Value = 40000
$ws.Range("L62").Value = 39999
$ws.Range("M62").Value = -39376
$ws.Range("N62").Value = -41247
$ws.Range("H65").Value = 39999.332
$ws.Range("I65").Value = 40000
$ws.Range("J65").Value = 39999
$ws.Range("K65").Value = 120000
$ws.Range("L65").Value = 119997
$ws.Range("M65").Value = -116880
$ws.Range("N65").Value = -126237
$ws.Range("H82").Value = 2670.2415
$ws.Range("I82").Value = 2717.35
$ws.Range("J82").Value = 2565.5557
$ws.Range("K82").Value = 2717.35
$ws.Range("L82").Value = 2565.5557
$ws.Range("M82").Value = -2356.35
$ws.Range("N82").Value = -3287.5557
$ws.Range("H85").Value = 2670.2415
$ws.Range("I85").Value = 2717.35
$ws.Range("J85").Value = 2565.5557
$ws.Range("K85").Value = 2717.35
$ws.Range("L85").Value = 2565.5557
$ws.Range("M85").Value = -1469.35
$ws.Range("N85").Value = -5061.5557
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524
$ws.Range("H132").Value = 1918769.2
$ws.Range("I132").Value = 2150.1667
$ws.Range("K132").Value = 6450.500100000001
$ws.Range("M132").Value = -3920.500100000001
$ws.Range("H136").Value = 20387.818
$ws.Range("I136").Value = 20866.545
$ws.Range("K136").Value = 62599.63499999999
$ws.Range("M136").Value = -60049.63499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 22000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26020
$ws.Range("H81").Value = 4183.4614
$ws.Range("I81").Value = 5226.1113
$ws.Range("J81").Value = 1837.5
$ws.Range("K81").Value = 10452.2226
$ws.Range("L81").Value = 3675
$ws.Range("M81").Value = -9391.222599999999
$ws.Range("N81").Value = -5797
$ws.Range("H84").Value = 4183.4614
$ws.Range("I84").Value = 5226.1113
$ws.Range("J84").Value = 1837.5
$ws.Range("K84").Value = 52261.113
$ws.Range("L84").Value = 18375
$ws.Range("M84").Value = -46957.113
$ws.Range("N84").Value = -28983
$ws.Range("H107").Value = 2012.375
$ws.Range("I107").Value = 2349.8333
$ws.Range("K107").Value = 7049.499899999999
$ws.Range("M107").Value = -5129.499899999999
$ws.Range("H132").Value = 10506.439
$ws.Range("I132").Value = 3493.6843
$ws.Range("K132").Value = 10481.0529
$ws.Range("M132").Value = -7951.052899999999
